$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three data rows for FERNANDO (004563237), LEONARDO (005076418)
# and CARLOS (005696533) - rows 3 through 5 of the "Export" sheet. Deleting
# the entire rows shifts everything below them up by three positions.
$ws.Range("A3:A5").EntireRow.Delete()

# The row that used to hold CAIO (004512434 / 3961.89) is now at row 5
# after the shift above; overwrite it with TATIANA's record.
$ws.Range("A5").Value = "'005348011"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "TATIANA"
$ws.Range("C5").Value = 2000
